$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assignment")

# ---------------------------------------------------------------------------
# 1. Insert two new columns (B:C) before the old "B" (subroutine depth) column.
#    This shifts old B -> D, old C -> E, old D -> F, old E -> G, old F -> H,
#    and Excel auto-adjusts every formula reference accordingly.
# ---------------------------------------------------------------------------
$ws.Range("B:C").Insert()

# ---------------------------------------------------------------------------
# 2. New column B: byte "Address" values (old local-index * 4), with its own
#    anchor + shared-style formulas, mirroring the old B(now D) column shape.
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "Address"

$ws.Range("B7").Formula = "=B8+4"
$ws.Range("B8:B20").FormulaR1C1 = "=R[1]C+4"
$ws.Range("B21").Formula = "=B22+4"
$ws.Range("B22").Value = 0

$ws.Range("B28").Formula = "=B29+4"
$ws.Range("B29:B33").FormulaR1C1 = "=R[1]C+4"
$ws.Range("B34").Formula = "=B35+4"
$ws.Range("B35").Value = 128

# ---------------------------------------------------------------------------
# 3. Bold header cell next to A4 ("Subroutine stack" / "Exception stack").
# ---------------------------------------------------------------------------
$ws.Range("B4").Font.Bold = $true
$ws.Range("B24").Font.Bold = $true

# ---------------------------------------------------------------------------
# 4. New columns J:L - hex "Address" lookup table mirrored alongside the
#    existing two tables (rows 7-22 and 28-35).
# ---------------------------------------------------------------------------
$ws.Range("J3").Value = "Address"
$ws.Range("J3:L3").Font.Bold = $true
$ws.Range("J3:L3").HorizontalAlignment = -4152

$ws.Range("J1:K2").HorizontalAlignment = -4152
$ws.Range("J4:K4").Font.Bold = $true
$ws.Range("J4:K4").HorizontalAlignment = -4152
$ws.Range("J5:K6").HorizontalAlignment = -4152

$ws.Range("J7:J22").FormulaR1C1 = "=DEC2HEX(RC[-8],2)"
$ws.Range("K7:K22").HorizontalAlignment = -4152

$ws.Range("J23:K23").HorizontalAlignment = -4152
$ws.Range("J24:K24").Font.Bold = $true
$ws.Range("J24:K24").HorizontalAlignment = -4152
$ws.Range("J25:K27").HorizontalAlignment = -4152

$ws.Range("J28:J35").FormulaR1C1 = "=DEC2HEX(RC[-8],2)"
$ws.Range("K28:K35").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# 5. Column widths / visibility.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 8.14
$ws.Columns.Item(2).Hidden = $true
$ws.Columns.Item(3).ColumnWidth = 5.14
$ws.Range("J1:K1").ColumnWidth = 11.28515625
$ws.Range("L1").ColumnWidth = 21
$ws.Range("M1").ColumnWidth = 35.42578125

# ---------------------------------------------------------------------------
# 6. Sheet view: select column K (matches the commit's saved selection) and
#    drop the frozen topLeftCell scroll position.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Columns.Item(11).Select()
